$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Step 0: in-place edits of the two rich-text header strings ---
# "Volume 30   Number  22" -> "...  23"  (only the trailing run changes)
$volRun = $ws.Range("A8").Characters(21, 2)
$volRun.Text = "23"

# "Report Covering the Week  5/29/2023  Through  6/4/2023"
#   -> "...  6/5/2023  Through  6/11/2023"
# Edit the right-most run first so the left run's character offsets stay valid.
$throughRun = $ws.Range("C9").Characters(47, 8)
$throughRun.Text = "6/11/2023"
$weekRun = $ws.Range("C9").Characters(27, 9)
$weekRun.Text = "6/5/2023"

# --- Step 1: fix cells whose literal/number vs text representation changed ---
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("D16").Copy($ws.Range("C16"))
$ws.Range("F20").Copy($ws.Range("D20"))
$ws.Range("H20").Copy($ws.Range("E20"))
$ws.Range("C22").Copy($ws.Range("D22"))
$ws.Range("N22").Copy($ws.Range("E22"))
$ws.Range("C22").Copy($ws.Range("C27"))
$ws.Range("I28").Copy($ws.Range("D28"))
$ws.Range("K28").Copy($ws.Range("E28"))
$ws.Range("I28").Copy($ws.Range("G28"))
$ws.Range("K28").Copy($ws.Range("H28"))
$ws.Range("I29").Copy($ws.Range("D29"))
$ws.Range("K29").Copy($ws.Range("E29"))
$ws.Range("I29").Copy($ws.Range("G29"))
$ws.Range("K29").Copy($ws.Range("H29"))

# --- Step 2: write the new numeric literal values ---
$ws.Range("N14").Value = -80
$ws.Range("M15").Value = -40
$ws.Range("N15").Value = -80
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -35.294117647058
$ws.Range("I16").Value = 73
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = -7.59493670886
$ws.Range("L16").Value = 7.35294117647
$ws.Range("M16").Value = -20.652173913043
$ws.Range("N16").Value = -79.494382022471
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 97
$ws.Range("K17").Value = 13.40206185567
$ws.Range("L17").Value = 10
$ws.Range("M17").Value = 29.411764705882
$ws.Range("N17").Value = -62.837837837837
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 7.142857142857
$ws.Range("I18").Value = 58
$ws.Range("J18").Value = 83
$ws.Range("K18").Value = -30.12048192771
$ws.Range("L18").Value = -9.375
$ws.Range("M18").Value = 11.538461538461
$ws.Range("N18").Value = -85.5
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -10
$ws.Range("I19").Value = 155
$ws.Range("J19").Value = 154
$ws.Range("K19").Value = 0.64935064935
$ws.Range("L19").Value = 6.896551724137
$ws.Range("M19").Value = 20.155038759689
$ws.Range("N19").Value = -11.428571428571
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = 13.043478260869
$ws.Range("M20").Value = 225
$ws.Range("N20").Value = -61.194029850746
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 66.666666666666
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = -5
$ws.Range("I21").Value = 428
$ws.Range("J21").Value = 440
$ws.Range("K21").Value = -2.727272727272
$ws.Range("L21").Value = 4.901960784313
$ws.Range("M21").Value = 14.745308310992
$ws.Range("N21").Value = -67.673716012084
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("L22").Value = -25
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -54.545454545454
$ws.Range("I23").Value = 36
$ws.Range("J23").Value = 46
$ws.Range("K23").Value = -21.739130434782
$ws.Range("L23").Value = 5.882352941176
$ws.Range("M23").Value = 157.142857142857
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 102
$ws.Range("G24").Value = 100
$ws.Range("H24").Value = 2
$ws.Range("I24").Value = 551
$ws.Range("J24").Value = 499
$ws.Range("K24").Value = 10.420841683366
$ws.Range("L24").Value = 20.568927789934
$ws.Range("M24").Value = 30.568720379146
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = -3.333333333333
$ws.Range("I25").Value = 189
$ws.Range("J25").Value = 191
$ws.Range("K25").Value = -1.047120418848
$ws.Range("L25").Value = 5
$ws.Range("M25").Value = -15.625
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 7
$ws.Range("K26").Value = 40
$ws.Range("L26").Value = -22.222222222222
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = -37.037037037037
$ws.Range("L27").Value = 0
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = -14.285714285714
$ws.Range("L28").Value = -50
$ws.Range("M28").Value = -40
$ws.Range("N28").Value = -86.363636363636
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -28.571428571428
$ws.Range("L29").Value = -58.333333333333
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -87.179487179487
